$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "52.022.41"
Set-TextCell "E2" "  +0.27%  "
Set-TextCell "D3" "2.981.94"
Set-TextCell "E3" "  +1.38%  "
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "354.54"
Set-TextCell "E5" "  +0.39%  "
Set-TextCell "D6" "107.56"
Set-TextCell "E6" "  -4.50%  "
Set-TextCell "D7" "0.562"
Set-TextCell "E7" "  +0.07%  "
Set-TextCell "E8" "  -0.04%  "
Set-TextCell "D9" "0.615"
Set-TextCell "E9" "  -2.03%  "
Set-TextCell "D10" "38.26"
Set-TextCell "E10" "  -3.00%  "
Set-TextCell "D12" "0.0857"
Set-TextCell "E12" "  -4.00%  "
Set-TextCell "D13" "19.27"
Set-TextCell "E13" "  -3.57%  "
Set-TextCell "D14" "3.448.43"
Set-TextCell "E14" "  +1.29%  "
Set-TextCell "E15" "  -3.14%  "
Set-TextCell "D16" "2.974.52"
Set-TextCell "E16" "  +1.16%  "
Set-TextCell "D17" "0.998"
Set-TextCell "E17" "  +0.65%  "
Set-TextCell "D18" "52.091.47"
Set-TextCell "E18" "  +0.25%  "
Set-TextCell "E19" "  +4.61%  "
Set-TextCell "E20" "  -2.19%  "
Set-TextCell "D21" "13.62"
Set-TextCell "E21" "  -5.86%  "
Set-TextCell "D22" "0.0₃0974"
Set-TextCell "E22" "  -1.79%  "
Set-TextCell "D23" "69.51"
Set-TextCell "E23" "  -2.50%  "
Set-TextCell "D24" "263.43"
Set-TextCell "E24" "  -2.50%  "
Set-TextCell "D25" "2.73"
Set-TextCell "E25" "  -2.15%  "
Set-TextCell "D26" "0.179"
Set-TextCell "E26" "  +0.25%  "
Set-TextCell "D27" "26.89"
Set-TextCell "E27" "  -0.29%  "
Set-TextCell "D28" "7.56"
Set-TextCell "E28" "  +2.40%  "
Set-TextCell "E29" "  -0.24%  "
Set-TextCell "D30" "0.109"
Set-TextCell "E30" "  +0.95%  "
Set-TextCell "D31" "10.30"
Set-TextCell "E31" "  -3.16%  "
Set-TextCell "D32" "6.09"
Set-TextCell "E32" "  -1.98%  "
Set-TextCell "D33" "36.29"
Set-TextCell "E33" "  -3.33%  "
Set-TextCell "E34" "  -3.43%  "
Set-TextCell "D35" "50.67"
Set-TextCell "E35" "  -4.08%  "
Set-TextCell "D36" "0.0447"
Set-TextCell "E36" "  -1.12%  "
Set-TextCell "D37" "0.998"
Set-TextCell "E37" "  +0.00%  "
Set-TextCell "D38" "3.19"
Set-TextCell "E38" "  -3.88%  "
Set-TextCell "D39" "17.90"
Set-TextCell "E39" "  -5.02%  "
Set-TextCell "E40" "  -3.66%  "
Set-TextCell "D41" "2.70"
Set-TextCell "E41" "  +0.11%  "
Set-TextCell "E42" "  -0.38%  "
Set-TextCell "D43" "22.76"
Set-TextCell "E43" "  -2.39%  "
Set-TextCell "D44" "121.95"
Set-TextCell "E44" "  +8.70%  "
Set-TextCell "E45" "  -3.04%  "
Set-TextCell "D46" "2.123.11"
Set-TextCell "E46" "  -2.13%  "
Set-TextCell "E47" "  -5.12%  "
Set-TextCell "D48" "2.36"
Set-TextCell "E48" "  -6.73%  "
Set-TextCell "D49" "0.241"
Set-TextCell "E49" "  -1.84%  "
Set-TextCell "D50" "0.0333"
Set-TextCell "E50" "  -2.94%  "
Set-TextCell "D51" "0.928"
Set-TextCell "E51" "  -1.01%  "

Write-Host "Updated cryptos list"
